$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed LR-pair TPM values (re-run of the NATMI scoring pipeline), now
# also covering the MuSCs target cluster in addition to FAPs for each
# sending cluster (Sending cluster | Ligand | Receptor | Target cluster |
# ... 16 numeric metric columns).
$rows = @()
$rows += ,@("ECs","Omg","Lingo1","FAPs",2,1,1.535573,3.071146,0.07799861199039625,0.05613057672622927,1,0.3333333333333333,0.096568,0.289704,0.3543662778110081,0.4515435187176778,0.148287213464,0.8897232807840001,0.02764007780546178,0.02534539812261415)
$rows += ,@("ECs","Omg","Lingo1","MuSCs",2,1,1.535573,3.071146,0.07799861199039625,0.05613057672622927,1,0.5,0.175941,0.351882,0.6456337221889918,0.5484564812823222,0.270170249193,1.080680996772,0.05035853418493446,0.03078517860361511)
$rows += ,@("FAPs","Omg","Lingo1","FAPs",3,1,1.766820333333333,5.300461,0.08974469701955506,0.09687521623683339,1,0.3333333333333333,0.096568,0.289704,0.3543662778110081,0.4515435187176778,0.1706183059493333,1.535564753544,0.0318024942360964,0.04374337601611566)
$rows += ,@("FAPs","Omg","Lingo1","MuSCs",3,1,1.766820333333333,5.300461,0.08974469701955506,0.09687521623683339,1,0.5,0.175941,0.351882,0.6456337221889918,0.5484564812823222,0.310856136267,1.865136817602,0.05794220278345866,0.05313184022071772)
$rows += ,@("Inflammatory-Mac","Omg","Lingo1","FAPs",3,1,4.075004666666667,12.225014,0.2069876898424154,0.2234335607314751,1,0.3333333333333333,0.096568,0.289704,0.3543662778110081,0.4515435187176778,0.3935150506506667,3.541635455856,0.07334945720215615,0.1008899762123102)
$rows += ,@("Inflammatory-Mac","Omg","Lingo1","MuSCs",3,1,4.075004666666667,12.225014,0.2069876898424154,0.2234335607314751,1,0.5,0.175941,0.351882,0.6456337221889918,0.5484564812823222,0.7169603960579999,4.301762376348,0.1336382326402592,0.1225435845191649)
$rows += ,@("MuSCs","Omg","Lingo1","FAPs",2,1,2.811664,5.623328,0.1428169741089258,0.1027761766326815,1,0.3333333333333333,0.096568,0.289704,0.3543662778110081,0.4515435187176778,0.271516769152,1.629100614912,0.05060951952321114,0.04640791643707059)
$rows += ,@("MuSCs","Omg","Lingo1","MuSCs",2,1,2.811664,5.623328,0.1428169741089258,0.1027761766326815,1,0.5,0.175941,0.351882,0.6456337221889918,0.5484564812823222,0.494686975824,1.978747903296,0.09220745458571462,0.05636826019561093)
$rows += ,@("Neutrophils","Omg","Lingo1","FAPs",3,1,7.050357000000001,21.151071,0.3581191255881512,0.3865728993696239,1,0.3333333333333333,0.096568,0.289704,0.3543662778110081,0.4515435187176778,0.6808388747760001,6.127549872984001,0.1269053415476061,0.1745544872222547)
$rows += ,@("Neutrophils","Omg","Lingo1","MuSCs",3,1,7.050357000000001,21.151071,0.3581191255881512,0.3865728993696239,1,0.5,0.175941,0.351882,0.6456337221889918,0.5484564812823222,1.240446860937,7.442681165622,0.2312137840405451,0.2120184121473691)
$rows += ,@("Resolving-Mac","Omg","Lingo1","FAPs",3,1,2.447764666666667,7.343294,0.1243329014505562,0.1342115703031569,1,0.3333333333333333,0.096568,0.289704,0.3543662778110081,0.4515435187176778,0.2363757383306667,2.127381644976,0.04405938749647649,0.06060236470731244)
$rows += ,@("Resolving-Mac","Omg","Lingo1","MuSCs",3,1,2.447764666666667,7.343294,0.1243329014505562,0.1342115703031569,1,0.5,0.175941,0.351882,0.6456337221889918,0.5484564812823222,0.4306621632179999,2.583972979308,0.08027351395407971,0.07360920559584441)

$nRows = $rows.Count
$nCols = 20
$arr = New-Object "object[,]" $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
  for ($c = 0; $c -lt $nCols; $c++) {
    $arr[$r,$c] = $rows[$r][$c]
  }
}

$target = $ws.Range("A2").Resize($nRows, $nCols)
$target.Value2 = $arr
